# Applies the coin-price/ranking refresh described in the commit
# "Updated symbol list on Sun Dec 18 08:35:48 UTC 2022 with GitHub Actions".
#
# All edited cells on the sheet are plain text cells (Coin / Link / Price /
# Volume(1h) columns stored as text, not numbers). Column D ("Price") holds
# numeric-looking text such as "247.43", so a bare assignment would make Excel
# auto-convert it to a Number cell; prefixing the literal with a single quote
# (the normal Excel "force text" input convention) keeps it a Text cell, just
# like the other text columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''247.43'

$ws.Range("D4").Value = '''5.551'

$ws.Range("D5").Value = '''0.05632'

$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '''6.480'
$ws.Range("E6").Value = '5KuCoinTokenKCS'

$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.074'
$ws.Range("E7").Value = '6FTXTokenFTT'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.8019'
$ws.Range("E8").Value = '7MXTokenMX'

$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1427'
$ws.Range("E9").Value = '8WazirXWRX'

$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '''0.07329'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.03190'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.02991'
$ws.Range("E12").Value = '11BitrueCoinBTR'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09263'
$ws.Range("E13").Value = '12BitMartTokenBMX'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001664'
$ws.Range("E14").Value = '13BitForexTokenBF'

$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").Value = '''0.04675'
$ws.Range("E15").Value = '14CoinExTokenCET'

$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = '''0.0005946'
$ws.Range("E16").Value = '15OneONE'

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.006274'
$ws.Range("E17").Value = '16TigerCashTCH'

$ws.Range("B18").Value = 'BitKan'
$ws.Range("C18").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D18").Value = '''0.001051'
$ws.Range("E18").Value = '17BitKanKAN'

$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = '''0.003833'
$ws.Range("E19").Value = '18HotbitTokenHTB'

$ws.Range("B20").Value = 'NitroEx'
$ws.Range("C20").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D20").Value = '''0.0001502'
$ws.Range("E20").Value = '19NitroExNTX'

$ws.Range("B21").Value = 'UpBots'
$ws.Range("C21").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D21").Value = '''0.0004604'
$ws.Range("E21").Value = '20UpBotsUBXT'

$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").Value = '''3.982'
$ws.Range("E22").Value = '21LEOLEO'

$ws.Range("B23").Value = 'GateToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D23").Value = '''3.396'
$ws.Range("E23").Value = '22GateTokenGT'

$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.113'
$ws.Range("E24").Value = '23BTSETokenBTSE'

$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '''0.3311'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'

$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '''0.1292'
$ws.Range("E26").Value = '25ProBitTokenPROB'

$ws.Range("B27").Value = 'MCDex'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D27").Value = '''2.588'
$ws.Range("E27").Value = '26MCDexMCB'

$ws.Range("D40").Value = '''0.04187'

$ws.Range("D41").Value = '''0.007027'

$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.003503'
$ws.Range("E42").Value = '41CEJICEJIBestin24h'

$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = '''0.1045'
$ws.Range("E43").Value = '42BKEXTokenBKK'

$ws.Range("D44").Value = '''0.008705'

$ws.Range("D45").Value = '''0.00005648'

$ws.Range("D46").Value = '''0.00000000751'

$ws.Range("D47").Value = '''0.6806'

$ws.Range("D48").Value = '''0.02750'

$ws.Range("D49").Value = '''0.00002102'

$ws.Range("D50").Value = '''0.01011'
